# Generate Report for Handoff
#
# The localization-status report re-evaluates the handoff queue. Three
# records (f2031f0a, 4a50c043, 88ef0b03) get re-ordered: the two rows that
# were already "Ready for handoff" move up (rows 7 & 8), and the row that
# was still "In Translation" drops to the bottom (row 9) with a refreshed
# "latest" timestamp, on all three sheets (Overview, zh-cn, de-de).
#
# This script overwrites the cell contents of rows 7-9 directly (rather
# than doing a physical row move), and separately fixes up the hyperlink
# display text so it mirrors the new cell text while the hyperlink
# targets (relationship ids / URLs) stay exactly as they were.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A-G, rows 7-9
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 7 <- old row 8 data
$ws.Range("A7").Value = "4a50c043-5595-44a0-a0f4-8ba6c0c90adf.md"
$ws.Range("B7").Value = "e2e\4a50c043-5595-44a0-a0f4-8ba6c0c90adf.md"
$ws.Range("C7").Value = ".md"
$ws.Range("E7").Value = "Ready for handoff"
$ws.Range("F7").Value = "Ready for handoff"
$ws.Range("G7").Value = "2016-09-06 17:16:21"

# Row 8 <- old row 9 data
$ws.Range("A8").Value = "88ef0b03-2d6b-40ae-8a32-a4f3ed3f9ade.md"
$ws.Range("B8").Value = "e2e\88ef0b03-2d6b-40ae-8a32-a4f3ed3f9ade.md"
$ws.Range("C8").Value = ".md"
$ws.Range("E8").Value = "Ready for handoff"
$ws.Range("F8").Value = "Ready for handoff"
$ws.Range("G8").Value = "2016-09-06 17:13:37"

# Row 9 <- old row 7 data, with a refreshed "Latest HO Xliff Generate Date"
$ws.Range("A9").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$ws.Range("B9").Value = "e2e\f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$ws.Range("C9").Value = ".md"
$ws.Range("E9").Value = "In Translation"
$ws.Range("F9").Value = "In Translation"
$ws.Range("G9").Value = "2016-09-06 17:19:36"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$7') { $hl.TextToDisplay = "e2e\4a50c043-5595-44a0-a0f4-8ba6c0c90adf.md" }
    elseif ($addr -eq '$B$8') { $hl.TextToDisplay = "e2e\88ef0b03-2d6b-40ae-8a32-a4f3ed3f9ade.md" }
    elseif ($addr -eq '$B$9') { $hl.TextToDisplay = "e2e\f2031f0a-2b1a-493f-893f-f63f1f1858ce.md" }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A-P, rows 7-9
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 7 <- old row 8 data
$ws.Range("A7").Value = "4a50c043-5595-44a0-a0f4-8ba6c0c90adf.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("G7").Value = "4a50c043-5595-44a0-a0f4-8ba6c0c90adf.4c439598832228233399d243f01db6f80e7586a9.zh-cn.xlf"
$ws.Range("H7").Value = "2016-09-06 17:16:14"

# Row 8 <- old row 9 data
$ws.Range("A8").Value = "88ef0b03-2d6b-40ae-8a32-a4f3ed3f9ade.md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("G8").Value = "88ef0b03-2d6b-40ae-8a32-a4f3ed3f9ade.abeeacf69803d2d65b20d8be402c4f1b04f0e525.zh-cn.xlf"
$ws.Range("H8").Value = "2016-09-06 17:13:31"

# Row 9 <- old row 7 data, with a refreshed "Latest Handoff Datetime"
$ws.Range("A9").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$ws.Range("C9").Value = "In Translation"
$ws.Range("G9").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.80888a3f371db147e2f85912f9532ae6b5ca5e8a.zh-cn.xlf"
$ws.Range("H9").Value = "2016-09-06 17:19:31"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$7') { $hl.TextToDisplay = "4a50c043-5595-44a0-a0f4-8ba6c0c90adf.md" }
    elseif ($addr -eq '$A$8') { $hl.TextToDisplay = "88ef0b03-2d6b-40ae-8a32-a4f3ed3f9ade.md" }
    elseif ($addr -eq '$A$9') { $hl.TextToDisplay = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md" }
}

# ---------------------------------------------------------------------
# Sheet "de-de": columns A-P, rows 7-9
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 7 <- old row 8 data
$ws.Range("A7").Value = "4a50c043-5595-44a0-a0f4-8ba6c0c90adf.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("G7").Value = "4a50c043-5595-44a0-a0f4-8ba6c0c90adf.4c439598832228233399d243f01db6f80e7586a9.de-de.xlf"
$ws.Range("H7").Value = "2016-09-06 17:16:21"

# Row 8 <- old row 9 data
$ws.Range("A8").Value = "88ef0b03-2d6b-40ae-8a32-a4f3ed3f9ade.md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("G8").Value = "88ef0b03-2d6b-40ae-8a32-a4f3ed3f9ade.abeeacf69803d2d65b20d8be402c4f1b04f0e525.de-de.xlf"
$ws.Range("H8").Value = "2016-09-06 17:13:37"

# Row 9 <- old row 7 data, with a refreshed "Latest Handoff Datetime"
$ws.Range("A9").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md"
$ws.Range("C9").Value = "In Translation"
$ws.Range("G9").Value = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.80888a3f371db147e2f85912f9532ae6b5ca5e8a.de-de.xlf"
$ws.Range("H9").Value = "2016-09-06 17:19:36"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$7') { $hl.TextToDisplay = "4a50c043-5595-44a0-a0f4-8ba6c0c90adf.md" }
    elseif ($addr -eq '$A$8') { $hl.TextToDisplay = "88ef0b03-2d6b-40ae-8a32-a4f3ed3f9ade.md" }
    elseif ($addr -eq '$A$9') { $hl.TextToDisplay = "f2031f0a-2b1a-493f-893f-f63f1f1858ce.md" }
}
